$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 25 de Abril de 2020 a las 00:52'
$ws.Range("A4").Value = 'Estados Unidos'
$ws.Range("B4").Value = 922293
$ws.Range("C4").Value = 35851
$ws.Range("D4").Value = 93283
$ws.Range("E4").Value = 776949
$ws.Range("F4").Value = 14946
$ws.Range("G4").Value = 1827
$ws.Range("H4").Value = 52061

$ws.Range("A43").Value = 'Noruega'
$ws.Range("B43").Value = 7463
$ws.Range("C43").Value = 62
$ws.Range("D43").Value = 32
$ws.Range("E43").Value = 7232
$ws.Range("F43").Value = 53
$ws.Range("G43").Value = 5
$ws.Range("H43").Value = 199

$ws.Range("A85").Value = 'Costa de Marfil'
$ws.Range("B85").Value = 1077
$ws.Range("C85").Value = 73
$ws.Range("D85").Value = 419
$ws.Range("E85").Value = 644
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 14

$ws.Range("A86").Value = 'Hong Kong'
$ws.Range("B86").Value = 1036
$ws.Range("C86").Value = 0
$ws.Range("D86").Value = 699
$ws.Range("E86").Value = 333
$ws.Range("F86").Value = 9
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 4

$ws.Range("A89").Value = 'Guinea'
$ws.Range("B89").Value = 954
$ws.Range("C89").Value = 92
$ws.Range("D89").Value = 191
$ws.Range("E89").Value = 757
$ws.Range("F89").Value = 0
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 6

$ws.Range("A90").Value = 'Tunez'
$ws.Range("B90").Value = 922
$ws.Range("C90").Value = 4
$ws.Range("D90").Value = 194
$ws.Range("E90").Value = 690
$ws.Range("F90").Value = 20
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 38

$ws.Range("A98").Value = 'Niger'
$ws.Range("B98").Value = 681
$ws.Range("C98").Value = 10
$ws.Range("D98").Value = 289
$ws.Range("E98").Value = 368
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 24

$ws.Range("A99").Value = 'Albania'
$ws.Range("B99").Value = 678
$ws.Range("C99").Value = 15
$ws.Range("D99").Value = 394
$ws.Range("E99").Value = 257
$ws.Range("F99").Value = 4
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 27

$ws.Range("A112").Value = 'Reunion'
$ws.Range("B112").Value = 412
$ws.Range("C112").Value = 0
$ws.Range("D112").Value = 300
$ws.Range("E112").Value = 112
$ws.Range("F112").Value = 2
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 0

$ws.Range("A119").Value = 'Mali'
$ws.Range("B119").Value = 325
$ws.Range("C119").Value = 16
$ws.Range("D119").Value = 87
$ws.Range("E119").Value = 217
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 21

$ws.Range("A120").Value = 'Montenegro'
$ws.Range("B120").Value = 319
$ws.Range("C120").Value = 3
$ws.Range("D120").Value = 123
$ws.Range("E120").Value = 190
$ws.Range("F120").Value = 7
$ws.Range("G120").Value = 1
$ws.Range("H120").Value = 6

$ws.Range("A121").Value = 'Venezuela'
$ws.Range("B121").Value = 318
$ws.Range("C121").Value = 20
$ws.Range("D121").Value = 132
$ws.Range("E121").Value = 176
$ws.Range("F121").Value = 4
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 10

$ws.Range("A131").Value = 'Ruanda'
$ws.Range("B131").Value = 176
$ws.Range("C131").Value = 22
$ws.Range("D131").Value = 87
$ws.Range("E131").Value = 89
$ws.Range("F131").Value = 0
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 0

$ws.Range("A132").Value = 'Sudan'
$ws.Range("B132").Value = 174
$ws.Range("C132").Value = 12
$ws.Range("D132").Value = 14
$ws.Range("E132").Value = 144
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 3
$ws.Range("H132").Value = 16

$ws.Range("A133").Value = 'Gabon'
$ws.Range("B133").Value = 172
$ws.Range("C133").Value = 5
$ws.Range("D133").Value = 26
$ws.Range("E133").Value = 143
$ws.Range("F133").Value = 1
$ws.Range("G133").Value = 1
$ws.Range("H133").Value = 3

$ws.Range("A134").Value = 'Martinica'
$ws.Range("B134").Value = 170
$ws.Range("C134").Value = 6
$ws.Range("D134").Value = 77
$ws.Range("E134").Value = 79
$ws.Range("F134").Value = 6
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 14

$ws.Range("A139").Value = 'Maldivas'
$ws.Range("B139").Value = 129
$ws.Range("C139").Value = 21
$ws.Range("D139").Value = 16
$ws.Range("E139").Value = 113
$ws.Range("F139").Value = 2
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 0

$ws.Range("A140").Value = 'Madagascar'
$ws.Range("B140").Value = 122
$ws.Range("C140").Value = 1
$ws.Range("D140").Value = 61
$ws.Range("E140").Value = 61
$ws.Range("F140").Value = 1
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 0

$ws.Range("A141").Value = 'Camboya'
$ws.Range("B141").Value = 122
$ws.Range("C141").Value = 0
$ws.Range("D141").Value = 110
$ws.Range("E141").Value = 12
$ws.Range("F141").Value = 1
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 0

$ws.Range("A142").Value = 'Etiopia'
$ws.Range("B142").Value = 117
$ws.Range("C142").Value = 1
$ws.Range("D142").Value = 25
$ws.Range("E142").Value = 89
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 3

$ws.Range("A143").Value = 'Liberia'
$ws.Range("B143").Value = 117
$ws.Range("C143").Value = 16
$ws.Range("D143").Value = 25
$ws.Range("E143").Value = 84
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 8

$ws.Range("A157").Value = 'Bahamas'
$ws.Range("B157").Value = 73
$ws.Range("C157").Value = 1
$ws.Range("D157").Value = 15
$ws.Range("E157").Value = 47
$ws.Range("F157").Value = 1
$ws.Range("G157").Value = 0
$ws.Range("H157").Value = 11

$ws.Range("A158").Value = 'San Martin (Parte Holandesa)'
$ws.Range("B158").Value = 73
$ws.Range("C158").Value = 0
$ws.Range("D158").Value = 22
$ws.Range("E158").Value = 39
$ws.Range("F158").Value = 8
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 12

$ws.Range("A159").Value = 'Haiti'
$ws.Range("B159").Value = 72
$ws.Range("C159").Value = 0
$ws.Range("D159").Value = 2
$ws.Range("E159").Value = 65
$ws.Range("F159").Value = 0
$ws.Range("G159").Value = 0
$ws.Range("H159").Value = 5

$ws.Range("A162").Value = 'Libia'
$ws.Range("B162").Value = 61
$ws.Range("C162").Value = 1
$ws.Range("D162").Value = 18
$ws.Range("E162").Value = 41
$ws.Range("F162").Value = 0
$ws.Range("G162").Value = 0
$ws.Range("H162").Value = 2

